# PCB Design Checklist - complete schematic section for 2_Interface_JTAG board
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: "Components power supplies (dual supply, voltage level)." -> add comment
$ws.Range("F10").Value = "Power provided from debugger"

# Row 14: "RC circuit (10k - 100nF) o uC RESET pin." -> add comment
$ws.Range("F14").Value = "No IC"

# Row 19: "Differential pair rules." -> add comment
$ws.Range("F19").Value = "No high speed design"

# Row 21: "ESD, EFT, Surge Protection" -> move mark from NOK(D) to N/A(E), add comment
$ws.Range("D21").Value = ""
$ws.Range("E21").Value = "X"
$ws.Range("F21").Value = "Not needed for this specific board"

# Row 23: "Add Port Cross Reference (R+P+D)." -> move mark from NOK(D) to N/A(E), add comment
$ws.Range("D23").Value = ""
$ws.Range("E23").Value = "X"
$ws.Range("F23").Value = "Simple PCB design (1 sheet)"

# Update the sheet view selection to match the saved state (select header row)
$ws.Range("A1:F1").Select() | Out-Null
